$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 0.21075711429381644
$ws.Range("D2").Value = 0.21075711429381644
$ws.Range("E2").Value = 0.10134734964029923
$ws.Range("F2").Value = 0.0006587510490031676
$ws.Range("G2").Value = 0.7514
$ws.Range("C3").Value = 19.949309312641933
$ws.Range("D3").Value = 19.949309312641933
$ws.Range("E3").Value = 9.593078899211893
$ws.Range("F3").Value = 0.06235437641393587
$ws.Range("G3").Value = 0.0017
$ws.Range("C4").Value = 0.3187909077606328
$ws.Range("D4").Value = 0.3187909077606328
$ws.Range("E4").Value = 0.15329785520750586
$ws.Range("F4").Value = 0.000996425888652198
$ws.Range("G4").Value = 0.6966
$ws.Range("C5").Value = 299.4555315558221
$ws.Range("D5").Value = 2.079552302470987
$ws.Range("F5").Value = 0.9359904466484088
$ws.Range("C6").Value = 319.9343888905185
$ws.Range("C7").Value = 0.014985244785041773
$ws.Range("D7").Value = 0.014985244785041773
$ws.Range("E7").Value = 0.018466824508721714
$ws.Range("F7").Value = 0.00012765746886156176
$ws.Range("G7").Value = 0.8936
$ws.Range("C8").Value = 0.09824699122162449
$ws.Range("D8").Value = 0.09824699122162449
$ws.Range("E8").Value = 0.12107309366148633
$ws.Range("F8").Value = 0.0008369541106953426
$ws.Range("G8").Value = 0.7284
$ws.Range("C9").Value = 0.42166987486295987
$ws.Range("D9").Value = 0.42166987486295987
$ws.Range("E9").Value = 0.5196380634023269
$ws.Range("F9").Value = 0.0035921541284336706
$ws.Range("G9").Value = 0.4741
$ws.Range("C10").Value = 116.85145153282149
$ws.Range("D10").Value = 0.8114684134223714
$ws.Range("F10").Value = 0.9954432342920094
$ws.Range("C11").Value = 117.38635364369111
$ws.Range("C12").Value = 0.024329534532857605
$ws.Range("D12").Value = 0.024329534532857605
$ws.Range("E12").Value = 0.07103581191338554
$ws.Range("F12").Value = 0.00048590278506825305
$ws.Range("G12").Value = 0.7934
$ws.Range("C13").Value = 0.03956697516087168
$ws.Range("D13").Value = 0.03956697516087168
$ws.Range("E13").Value = 0.11552511215179229
$ws.Range("F13").Value = 0.000790220766510316
$ws.Range("G13").Value = 0.7376
$ws.Range("C14").Value = 0.6873581967692158
$ws.Range("D14").Value = 0.6873581967692158
$ws.Range("E14").Value = 2.0069043045965302
$ws.Range("F14").Value = 0.013727729221395254
$ws.Range("G14").Value = 0.1666
$ws.Range("C15").Value = 49.31953163290764
$ws.Range("D15").Value = 0.3424967474507475
$ws.Range("F15").Value = 0.9849961472270262
$ws.Range("C16").Value = 50.070786339370585
$ws.Range("C17").Value = 0.00018535900938485854
$ws.Range("D17").Value = 0.00018535900938485854
$ws.Range("E17").Value = 0.00024462077924020657
$ws.Range("F17").Value = 0.0000016791781192957105
$ws.Range("G17").Value = 0.9889
$ws.Range("C18").Value = 1.2000776516204612
$ws.Range("D18").Value = 1.2000776516204612
$ws.Range("E18").Value = 1.583758627446219
$ws.Range("F18").Value = 0.010871573713867023
$ws.Range("G18").Value = 0.2143
$ws.Range("C19").Value = 0.07188809784372463
$ws.Range("D19").Value = 0.07188809784372463
$ws.Range("E19").Value = 0.09487169019184793
$ws.Range("F19").Value = 0.0006512384876115563
$ws.Range("G19").Value = 0.7615
$ws.Range("C20").Value = 109.11459539260802
$ws.Range("D20").Value = 0.7577402457820002
$ws.Range("F20").Value = 0.9884755086204021
$ws.Range("C21").Value = 110.38674650108159
$ws.Range("C22").Value = 0.5151029099627621
$ws.Range("D22").Value = 0.5151029099627621
$ws.Range("E22").Value = 3.035545453750894
$ws.Range("F22").Value = 0.01979756498318349
$ws.Range("G22").Value = 0.085
$ws.Range("C23").Value = 0.6897191359769036
$ws.Range("D23").Value = 0.6897191359769036
$ws.Range("E23").Value = 4.06457379114989
$ws.Range("F23").Value = 0.026508798825529933
$ws.Range("G23").Value = 0.0456
$ws.Range("C24").Value = 0.3782591923860461
$ws.Range("D24").Value = 0.3782591923860461
$ws.Range("E24").Value = 2.229113735486283
$ws.Range("F24").Value = 0.014538087044180399
$ws.Range("G24").Value = 0.1359
$ws.Range("C25").Value = 24.43541701639917
$ws.Range("D25").Value = 0.16969039594721647
$ws.Range("F25").Value = 0.9391555491471062
$ws.Range("C26").Value = 26.018498254724882
$ws.Range("C27").Value = 0.024504662180262204
$ws.Range("D27").Value = 0.024504662180262204
$ws.Range("E27").Value = 0.05037321734871641
$ws.Range("F27").Value = 0.0003488594343584305
$ws.Range("G27").Value = 0.8171
$ws.Range("C28").Value = 0.05079873230225809
$ws.Range("D28").Value = 0.05079873230225809
$ws.Range("E28").Value = 0.10442484636095184
$ws.Range("F28").Value = 0.0007231936880715433
$ws.Range("G28").Value = 0.7444
$ws.Range("C29").Value = 0.11637409341924573
$ws.Range("D29").Value = 0.11637409341924573
$ws.Range("E29").Value = 0.23922539549593444
$ws.Range("F29").Value = 0.001656754135419744
$ws.Range("G29").Value = 0.6292
$ws.Range("C30").Value = 70.05054550178883
$ws.Range("D30").Value = 0.48646212154020024
$ws.Range("F30").Value = 0.9972711927421503
$ws.Range("C31").Value = 70.2422229896906
$ws.Range("C32").Value = 0.02198466568461665
$ws.Range("D32").Value = 0.02198466568461665
$ws.Range("E32").Value = 0.018768248912625986
$ws.Range("F32").Value = 0.00011972302791403514
$ws.Range("G32").Value = 0.8915
$ws.Range("C33").Value = 14.916104555775465
$ws.Range("D33").Value = 14.916104555775465
$ws.Range("E33").Value = 12.733837626898165
$ws.Range("F33").Value = 0.0812293999698778
$ws.Range("G33").Value = 0.0006
$ws.Range("C34").Value = 0.01323100014625389
$ws.Range("D34").Value = 0.01323100014625389
$ws.Range("E34").Value = 0.011295268605409949
$ws.Range("F34").Value = 0.00007205273996724773
$ws.Range("G34").Value = 0.9142
$ws.Range("C35").Value = 168.67806225945088
$ws.Range("D35").Value = 1.1713754323572978
$ws.Range("F35").Value = 0.918578824262241
$ws.Range("C36").Value = 183.6293824810572
$ws.Range("C37").Value = 0.07071000669817332
$ws.Range("D37").Value = 0.07071000669817332
$ws.Range("E37").Value = 0.2027099944885343
$ws.Range("F37").Value = 0.0013579457065837552
$ws.Range("G37").Value = 0.6929
$ws.Range("C38").Value = 1.703025948354592
$ws.Range("D38").Value = 1.703025948354592
$ws.Range("E38").Value = 4.882199800635975
$ws.Range("F38").Value = 0.03270565062509871
$ws.Range("G38").Value = 0.0528
$ws.Range("C39").Value = 0.06698842739827304
$ws.Range("D39").Value = 0.06698842739827304
$ws.Range("E39").Value = 0.19204104741020064
$ws.Range("F39").Value = 0.0012864748799215212
$ws.Range("G39").Value = 0.7341
$ws.Range("C40").Value = 50.230581823201064
$ws.Range("D40").Value = 0.3488234848833407
$ws.Range("F40").Value = 0.9646499287883961
$ws.Range("C41").Value = 52.0713062056521
